$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> (stage_life [H], dimension [I]) pairs to add for the "Adds briefs to main pipeline" edit.
$data = @(
    @(2, 'Adulthood and Elderly', 'Health'),
    @(3, 'School-aged Children', 'Education'),
    @(4, 'Prenatal and Early Childhood ', 'Health'),
    @(5, 'School-aged Children', 'Education'),
    @(6, 'All Ages', $null),
    @(7, 'All Ages', $null),
    @(8, 'All Ages', $null),
    @(9, 'Prenatal and Early Childhood ', 'Health'),
    @(10, 'School-aged Children', 'Education'),
    @(11, 'Adulthood and Elderly', 'Labor'),
    @(12, 'Adulthood and Elderly', 'Labor'),
    @(13, 'Adulthood and Elderly', 'Labor'),
    @(14, 'Adulthood and Elderly', 'Labor'),
    @(15, 'Adulthood and Elderly', 'Labor'),
    @(16, 'Adulthood and Elderly', 'Labor'),
    @(17, 'Youth', 'Labor'),
    @(18, 'Youth', 'Labor'),
    @(19, 'Youth', 'Labor'),
    @(20, 'Youth', 'Labor'),
    @(21, 'Youth', 'Labor'),
    @(22, 'Youth', 'Labor'),
    @(23, 'School-aged Children', 'Education'),
    @(24, 'All Ages', 'Education'),
    @(25, 'All Ages', 'Education'),
    @(26, 'Prenatal and Early Childhood ', 'Education'),
    @(27, 'School-aged Children', 'Education'),
    @(28, 'School-aged Children', 'Education'),
    @(29, 'School-aged Children', 'Education'),
    @(30, 'School-aged Children', 'Education'),
    @(31, 'Youth', 'Education'),
    @(32, 'All Ages', 'Health'),
    @(33, 'All Ages', 'Health'),
    @(34, 'School-aged Children', 'Education'),
    @(35, 'All Ages', 'Health'),
    @(36, 'All Ages', 'Health'),
    @(37, 'School-aged Children', 'Education'),
    @(38, 'All Ages', 'Health'),
    @(39, 'All Ages', 'Health'),
    @(40, 'School-aged Children', 'Education'),
    @(41, 'Prenatal and Early Childhood ', 'Health'),
    @(42, 'Prenatal and Early Childhood ', 'Health'),
    @(43, 'Prenatal and Early Childhood ', 'Health'),
    @(44, 'Prenatal and Early Childhood ', 'Health'),
    @(45, 'Prenatal and Early Childhood ', 'Health'),
    @(46, 'Youth', 'Health'),
    @(47, 'School-aged Children', 'Health'),
    @(48, 'Prenatal and Early Childhood ', 'Health'),
    @(49, 'Prenatal and Early Childhood ', 'Health'),
    @(50, 'Prenatal and Early Childhood ', 'Health'),
    @(51, 'Prenatal and Early Childhood ', 'Health'),
    @(52, 'Prenatal and Early Childhood ', 'Health'),
    @(53, 'Prenatal and Early Childhood ', 'Health'),
    @(54, 'Prenatal and Early Childhood ', 'Health'),
    @(55, 'Prenatal and Early Childhood ', 'Health'),
    @(56, 'Prenatal and Early Childhood ', 'Health'),
    @(57, 'Prenatal and Early Childhood ', 'Health'),
    @(58, 'Prenatal and Early Childhood ', 'Health'),
    @(59, 'Prenatal and Early Childhood ', 'Health'),
    @(60, 'Prenatal and Early Childhood ', 'Health'),
    @(61, 'Prenatal and Early Childhood ', 'Health'),
    @(62, 'Prenatal and Early Childhood ', 'Health'),
    @(63, 'School-aged Children', 'Education'),
    @(64, 'Youth', 'Health'),
    @(65, 'All Ages', 'Health'),
    @(66, 'Prenatal and Early Childhood ', 'Health'),
    @(67, 'All Ages', 'Health'),
    @(68, 'School-aged Children', 'Education'),
    @(69, 'Prenatal and Early Childhood ', 'Health'),
    @(70, 'Prenatal and Early Childhood ', 'Health'),
    @(71, 'School-aged Children', 'Education'),
    @(72, 'School-aged Children', 'Education'),
    @(73, 'School-aged Children', 'Education'),
    @(74, 'Youth', 'Education'),
)

foreach ($entry in $data) {
    $r = $entry[0]
    $h = $entry[1]
    $i = $entry[2]
    if ($null -ne $h) {
        $ws.Cells.Item($r, 8).Value = $h
    }
    if ($null -ne $i) {
        $ws.Cells.Item($r, 9).Value = $i
    }
}

# Column J: rank counter, 1 at J2 then incrementing by 1 down to J74 (each
# row's formula referencing the row directly above it, same as an Excel
# fill-down of "=J2+1").
$ws.Range("J2").Formula = "=1"
for ($r = 3; $r -le 74; $r++) {
    $prev = $r - 1
    $ws.Cells.Item($r, 10).Formula = "=J" + $prev + "+1"
}

# Final selection/view state left by the author after editing.
[void]$ws.Range("J3").Select()
